$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating-point precision of the existing A46 timestamp
$ws.Range("A46").Value = 44359.76822617477

# Append new row 47 with the latest retrieved data
$ws.Range("A47").Value = 44360.76849045554
$ws.Range("B47").Value = 76651
$ws.Range("C47").Value = 64418
$ws.Range("D47").Value = 3403
$ws.Range("E47").Value = 2083
$ws.Range("F47").Value = 1461
$ws.Range("G47").Value = 20242
$ws.Range("H47").Value = 1488
$ws.Range("I47").Value = 877
$ws.Range("J47").Value = 188

# Match the date-time style used by the rest of column A
$ws.Range("A46").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Range("A47").NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
